$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header values P1:Q1 (14, 15), copying the style of the existing header cell O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Cells.Item(1, 16).Value = 14   # P1
$ws.Cells.Item(1, 17).Value = 15   # Q1

# Update data values for columns B:I (rows 2-25) with recalculated results
# Row 2
$ws.Cells.Item(2, 2).Value = 24.44307359118704
$ws.Cells.Item(2, 3).Value = 18.84945595138836
$ws.Cells.Item(2, 4).Value = 10.02066305241336
$ws.Cells.Item(2, 5).Value = 29.10526702065993
$ws.Cells.Item(2, 6).Value = 63.03962927701833
$ws.Cells.Item(2, 7).Value = 2.074527469938671
$ws.Cells.Item(2, 8).Value = 3.089966362669392
$ws.Cells.Item(2, 9).Value = 3.037267292020123
# Row 3
$ws.Cells.Item(3, 2).Value = 22.75403469952843
$ws.Cells.Item(3, 3).Value = 17.53033396951019
$ws.Cells.Item(3, 4).Value = 9.477274155190834
$ws.Cells.Item(3, 5).Value = 27.02775996347093
$ws.Cells.Item(3, 6).Value = 59.07048019432784
$ws.Cells.Item(3, 7).Value = 2.084514762941586
$ws.Cells.Item(3, 8).Value = 2.711583836633365
$ws.Cells.Item(3, 9).Value = 2.666704839883651
# Row 4
$ws.Cells.Item(4, 2).Value = 21.65952378285824
$ws.Cells.Item(4, 3).Value = 16.68690558818472
$ws.Cells.Item(4, 4).Value = 9.126601439550928
$ws.Cells.Item(4, 5).Value = 25.68821643124625
$ws.Cells.Item(4, 6).Value = 56.50830951348823
$ws.Cells.Item(4, 7).Value = 2.090765124095499
$ws.Cells.Item(4, 8).Value = 2.473648940128331
$ws.Cells.Item(4, 9).Value = 2.509508808433992
# Row 5
$ws.Cells.Item(5, 2).Value = 21.19839618377675
$ws.Cells.Item(5, 3).Value = 16.35017261237931
$ws.Cells.Item(5, 4).Value = 8.962248306382884
$ws.Cells.Item(5, 5).Value = 25.1250067760214
$ws.Cells.Item(5, 6).Value = 55.35697436933109
$ws.Cells.Item(5, 7).Value = 2.093373239817404
$ws.Cells.Item(5, 8).Value = 2.374520949779008
$ws.Cells.Item(5, 9).Value = 2.613255113535315
# Row 6
$ws.Cells.Item(6, 2).Value = 21.12034009970148
$ws.Cells.Item(6, 3).Value = 16.31276175098225
$ws.Cells.Item(6, 4).Value = 8.913715193426158
$ws.Cells.Item(6, 5).Value = 25.02928931015711
$ws.Cells.Item(6, 6).Value = 55.07122712887607
$ws.Cells.Item(6, 7).Value = 2.093843404499031
$ws.Cells.Item(6, 8).Value = 2.357148475066212
$ws.Cells.Item(6, 9).Value = 2.635588437125068
# Row 7
$ws.Cells.Item(7, 2).Value = 21.65169548959532
$ws.Cells.Item(7, 3).Value = 16.73292013436325
$ws.Cells.Item(7, 4).Value = 9.067926423508622
$ws.Cells.Item(7, 5).Value = 25.67736305216642
$ws.Cells.Item(7, 6).Value = 56.24481937174074
$ws.Cells.Item(7, 7).Value = 2.090896024981767
$ws.Cells.Item(7, 8).Value = 2.470075046960208
$ws.Cells.Item(7, 9).Value = 2.524493452909088
# Row 8
$ws.Cells.Item(8, 2).Value = 23.870143067299
$ws.Cells.Item(8, 3).Value = 18.46235242138411
$ws.Cells.Item(8, 4).Value = 9.765799130928215
$ws.Cells.Item(8, 5).Value = 28.39744764467337
$ws.Cells.Item(8, 6).Value = 61.38995501941691
$ws.Cells.Item(8, 7).Value = 2.078075537225312
$ws.Cells.Item(8, 8).Value = 2.957617360055809
$ws.Cells.Item(8, 9).Value = 2.908925381665187
# Row 9
$ws.Cells.Item(9, 2).Value = 27.77635781228673
$ws.Cells.Item(9, 3).Value = 21.51752086636726
$ws.Cells.Item(9, 4).Value = 11.09435833615976
$ws.Cells.Item(9, 5).Value = 33.24354161310295
$ws.Cells.Item(9, 6).Value = 70.89697344570428
$ws.Cells.Item(9, 7).Value = 2.053570185398751
$ws.Cells.Item(9, 8).Value = 3.877518438033591
$ws.Cells.Item(9, 9).Value = 3.818003966742009
# Row 10
$ws.Cells.Item(10, 2).Value = 30.36896984993525
$ws.Cells.Item(10, 3).Value = 23.60291155827736
$ws.Cells.Item(10, 4).Value = 11.73454473692673
$ws.Cells.Item(10, 5).Value = 35.62103610224544
$ws.Cells.Item(10, 6).Value = 76.20990935265392
$ws.Cells.Item(10, 7).Value = 2.036709757866656
$ws.Cells.Item(10, 8).Value = 4.482358605421413
$ws.Cells.Item(10, 9).Value = 4.449622519871673
# Row 11
$ws.Cells.Item(11, 2).Value = 31.36121202343197
$ws.Cells.Item(11, 3).Value = 24.30366210427836
$ws.Cells.Item(11, 4).Value = 9.970562816033523
$ws.Cells.Item(11, 5).Value = 29.02521488677462
$ws.Cells.Item(11, 6).Value = 69.78476296936714
$ws.Cells.Item(11, 7).Value = 2.03570407727512
$ws.Cells.Item(11, 8).Value = 4.790747882991368
$ws.Cells.Item(11, 9).Value = 4.510909130059678
# Row 12
$ws.Cells.Item(12, 2).Value = 31.68585877720369
$ws.Cells.Item(12, 3).Value = 24.44128593523203
$ws.Cells.Item(12, 4).Value = 8.470649040072542
$ws.Cells.Item(12, 5).Value = 23.0426619106675
$ws.Cells.Item(12, 6).Value = 63.59460163662632
$ws.Cells.Item(12, 7).Value = 2.037548408334782
$ws.Cells.Item(12, 8).Value = 5.512775492837027
$ws.Cells.Item(12, 9).Value = 4.456295562001463
# Row 13
$ws.Cells.Item(13, 2).Value = 31.52803904234584
$ws.Cells.Item(13, 3).Value = 24.22003956114179
$ws.Cells.Item(13, 4).Value = 7.009883965657012
$ws.Cells.Item(13, 5).Value = 17.05336036085701
$ws.Cells.Item(13, 6).Value = 56.73539812628097
$ws.Cells.Item(13, 7).Value = 2.041798004963339
$ws.Cells.Item(13, 8).Value = 6.450362237943956
$ws.Cells.Item(13, 9).Value = 4.308155851271346
# Row 14
$ws.Cells.Item(14, 2).Value = 31.20074615630609
$ws.Cells.Item(14, 3).Value = 23.91840789243268
$ws.Cells.Item(14, 4).Value = 6.025607134096393
$ws.Cells.Item(14, 5).Value = 12.91878041738623
$ws.Cells.Item(14, 6).Value = 51.49898527829699
$ws.Cells.Item(14, 7).Value = 2.045765259257151
$ws.Cells.Item(14, 8).Value = 7.201276089694248
$ws.Cells.Item(14, 9).Value = 4.16614545595082
# Row 15
$ws.Cells.Item(15, 2).Value = 31.01509405687257
$ws.Cells.Item(15, 3).Value = 23.77999269105106
$ws.Cells.Item(15, 4).Value = 5.778124340762592
$ws.Cells.Item(15, 5).Value = 11.90505057293832
$ws.Cells.Item(15, 6).Value = 49.98084492291337
$ws.Cells.Item(15, 7).Value = 2.047349627796285
$ws.Cells.Item(15, 8).Value = 7.37155626470092
$ws.Cells.Item(15, 9).Value = 4.109105430901421
# Row 16
$ws.Cells.Item(16, 2).Value = 29.98938299365462
$ws.Cells.Item(16, 3).Value = 23.00271204771266
$ws.Cells.Item(16, 4).Value = 5.728103560500506
$ws.Cells.Item(16, 5).Value = 11.57793924659197
$ws.Cells.Item(16, 6).Value = 48.56665406277664
$ws.Cells.Item(16, 7).Value = 2.053498711515744
$ws.Cells.Item(16, 8).Value = 7.055752869198356
$ws.Cells.Item(16, 9).Value = 3.877767251533557
# Row 17
$ws.Cells.Item(17, 2).Value = 29.36448092993919
$ws.Cells.Item(17, 3).Value = 22.55793443121987
$ws.Cells.Item(17, 4).Value = 6.19513540705
$ws.Cells.Item(17, 5).Value = 13.59355848671587
$ws.Cells.Item(17, 6).Value = 50.37157133274903
$ws.Cells.Item(17, 7).Value = 2.056209579269872
$ws.Cells.Item(17, 8).Value = 6.353811635908204
$ws.Cells.Item(17, 9).Value = 3.772518732282649
# Row 18
$ws.Cells.Item(18, 2).Value = 29.0326828153557
$ws.Cells.Item(18, 3).Value = 22.33313476896589
$ws.Cells.Item(18, 4).Value = 7.282784649532969
$ws.Cells.Item(18, 5).Value = 18.16992540045136
$ws.Cells.Item(18, 6).Value = 55.3707338003543
$ws.Cells.Item(18, 7).Value = 2.055887887705688
$ws.Cells.Item(18, 8).Value = 5.333293727927019
$ws.Cells.Item(18, 9).Value = 3.773529732633144
# Row 19
$ws.Cells.Item(19, 2).Value = 28.97345199471054
$ws.Cells.Item(19, 3).Value = 22.4033883500889
$ws.Cells.Item(19, 4).Value = 8.790375362995041
$ws.Cells.Item(19, 5).Value = 24.53208797344972
$ws.Cells.Item(19, 6).Value = 62.1806771869171
$ws.Cells.Item(19, 7).Value = 2.052752511625638
$ws.Cells.Item(19, 8).Value = 4.411514987705797
$ws.Cells.Item(19, 9).Value = 3.879578676368524
# Row 20
$ws.Cells.Item(20, 2).Value = 29.69820850587216
$ws.Cells.Item(20, 3).Value = 23.17917654727976
$ws.Cells.Item(20, 4).Value = 11.41151704201456
$ws.Cells.Item(20, 5).Value = 34.95363493126855
$ws.Cells.Item(20, 6).Value = 74.17869503357267
$ws.Cells.Item(20, 7).Value = 2.041462083656163
$ws.Cells.Item(20, 8).Value = 4.314225926561374
$ws.Cells.Item(20, 9).Value = 4.278525854025554
# Row 21
$ws.Cells.Item(21, 2).Value = 31.63458094084674
$ws.Cells.Item(21, 3).Value = 24.75784781552005
$ws.Cells.Item(21, 4).Value = 12.27927442675885
$ws.Cells.Item(21, 5).Value = 38.09980501401512
$ws.Cells.Item(21, 6).Value = 79.78419766314934
$ws.Cells.Item(21, 7).Value = 2.027249542232134
$ws.Cells.Item(21, 8).Value = 4.857352354950227
$ws.Cells.Item(21, 9).Value = 4.805519619189777
# Row 22
$ws.Cells.Item(22, 2).Value = 32.85232260185697
$ws.Cells.Item(22, 3).Value = 25.7093163221756
$ws.Cells.Item(22, 4).Value = 12.77195792536501
$ws.Cells.Item(22, 5).Value = 39.65220401541557
$ws.Cells.Item(22, 6).Value = 83.08140148289172
$ws.Cells.Item(22, 7).Value = 2.018208876805866
$ws.Cells.Item(22, 8).Value = 5.187829086509634
$ws.Cells.Item(22, 9).Value = 5.140100061090803
# Row 23
$ws.Cells.Item(23, 2).Value = 32.20893271165753
$ws.Cells.Item(23, 3).Value = 25.16305894123705
$ws.Cells.Item(23, 4).Value = 12.56410863771151
$ws.Cells.Item(23, 5).Value = 38.83252030482961
$ws.Cells.Item(23, 6).Value = 81.55397815487034
$ws.Cells.Item(23, 7).Value = 2.022918405131028
$ws.Cells.Item(23, 8).Value = 5.014325094539056
$ws.Cells.Item(23, 9).Value = 4.96347957896903
# Row 24
$ws.Cells.Item(24, 2).Value = 29.68185813993916
$ws.Cells.Item(24, 3).Value = 23.11238178875517
$ws.Cells.Item(24, 4).Value = 11.66976099974378
$ws.Cells.Item(24, 5).Value = 35.63223609375809
$ws.Cells.Item(24, 6).Value = 75.24053863113274
$ws.Cells.Item(24, 7).Value = 2.04088703992831
$ws.Cells.Item(24, 8).Value = 4.351378196554124
$ws.Cells.Item(24, 9).Value = 4.293550467567776
# Row 25
$ws.Cells.Item(25, 2).Value = 26.76275660477796
$ws.Cells.Item(25, 3).Value = 20.79094544270774
$ws.Cells.Item(25, 4).Value = 10.65831984704363
$ws.Cells.Item(25, 5).Value = 31.97679880725834
$ws.Cells.Item(25, 6).Value = 68.04404505137312
$ws.Cells.Item(25, 7).Value = 2.060308407478158
$ws.Cells.Item(25, 8).Value = 3.628166126820007
$ws.Cells.Item(25, 9).Value = 3.571692568677724

# Populate new columns P and Q (all zeros) for data rows 2-25
$ws.Cells.Item(2, 16).Value = 0   # P2
$ws.Cells.Item(2, 17).Value = 0   # Q2
$ws.Cells.Item(3, 16).Value = 0   # P3
$ws.Cells.Item(3, 17).Value = 0   # Q3
$ws.Cells.Item(4, 16).Value = 0   # P4
$ws.Cells.Item(4, 17).Value = 0   # Q4
$ws.Cells.Item(5, 16).Value = 0   # P5
$ws.Cells.Item(5, 17).Value = 0   # Q5
$ws.Cells.Item(6, 16).Value = 0   # P6
$ws.Cells.Item(6, 17).Value = 0   # Q6
$ws.Cells.Item(7, 16).Value = 0   # P7
$ws.Cells.Item(7, 17).Value = 0   # Q7
$ws.Cells.Item(8, 16).Value = 0   # P8
$ws.Cells.Item(8, 17).Value = 0   # Q8
$ws.Cells.Item(9, 16).Value = 0   # P9
$ws.Cells.Item(9, 17).Value = 0   # Q9
$ws.Cells.Item(10, 16).Value = 0   # P10
$ws.Cells.Item(10, 17).Value = 0   # Q10
$ws.Cells.Item(11, 16).Value = 0   # P11
$ws.Cells.Item(11, 17).Value = 0   # Q11
$ws.Cells.Item(12, 16).Value = 0   # P12
$ws.Cells.Item(12, 17).Value = 0   # Q12
$ws.Cells.Item(13, 16).Value = 0   # P13
$ws.Cells.Item(13, 17).Value = 0   # Q13
$ws.Cells.Item(14, 16).Value = 0   # P14
$ws.Cells.Item(14, 17).Value = 0   # Q14
$ws.Cells.Item(15, 16).Value = 0   # P15
$ws.Cells.Item(15, 17).Value = 0   # Q15
$ws.Cells.Item(16, 16).Value = 0   # P16
$ws.Cells.Item(16, 17).Value = 0   # Q16
$ws.Cells.Item(17, 16).Value = 0   # P17
$ws.Cells.Item(17, 17).Value = 0   # Q17
$ws.Cells.Item(18, 16).Value = 0   # P18
$ws.Cells.Item(18, 17).Value = 0   # Q18
$ws.Cells.Item(19, 16).Value = 0   # P19
$ws.Cells.Item(19, 17).Value = 0   # Q19
$ws.Cells.Item(20, 16).Value = 0   # P20
$ws.Cells.Item(20, 17).Value = 0   # Q20
$ws.Cells.Item(21, 16).Value = 0   # P21
$ws.Cells.Item(21, 17).Value = 0   # Q21
$ws.Cells.Item(22, 16).Value = 0   # P22
$ws.Cells.Item(22, 17).Value = 0   # Q22
$ws.Cells.Item(23, 16).Value = 0   # P23
$ws.Cells.Item(23, 17).Value = 0   # Q23
$ws.Cells.Item(24, 16).Value = 0   # P24
$ws.Cells.Item(24, 17).Value = 0   # Q24
$ws.Cells.Item(25, 16).Value = 0   # P25
$ws.Cells.Item(25, 17).Value = 0   # Q25
